$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.866.85'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '1.639.66'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('D4').Value = '''0.9997'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '''309.26'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '''0.3868'
$ws.Range('E7').Value = '  -1.02%  '
$ws.Range('D8').Value = '''0.3824'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('D9').Value = '''50.39'
$ws.Range('E9').Value = '  -2.43%  '
$ws.Range('D10').Value = '''1.324'
$ws.Range('E10').Value = '  -3.70%  '
$ws.Range('D11').Value = '''0.9994'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').Value = '''0.08369'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').Value = '''23.76'
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('D14').Value = '''6.954'
$ws.Range('E14').Value = '  -4.12%  '
$ws.Range('D15').Value = '''7.774'
$ws.Range('E15').Value = '  -2.88%  '
$ws.Range('D16').Value = '''0.00001307'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').Value = '1.638.24'
$ws.Range('E17').Value = '  -1.44%  '
$ws.Range('D18').Value = '''93.48'
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('D19').Value = '''0.06940'
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('D20').Value = '''19.40'
$ws.Range('E20').Value = '  -2.97%  '
$ws.Range('D21').Value = '''6.869'
$ws.Range('E21').Value = '  -2.12%  '
$ws.Range('D22').Value = '''1.001'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = '''13.52'
$ws.Range('E23').Value = '  -1.73%  '
$ws.Range('D24').Value = '23.853.38'
$ws.Range('E24').Value = '  -0.70%  '
$ws.Range('D25').Value = '''2.433'
$ws.Range('E25').Value = '  -2.25%  '
$ws.Range('D26').Value = '''2.877'
$ws.Range('E26').Value = '  -9.06%  '
$ws.Range('D27').Value = '''21.84'
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('D28').Value = '''152.93'
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('D29').Value = '''5.468'
$ws.Range('E29').Value = '  +3.20%  '
$ws.Range('D30').Value = '''136.33'
$ws.Range('E30').Value = '  -2.79%  '
$ws.Range('D31').Value = '''7.797'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').Value = '''2.479'
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').Value = '1.819.68'
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').Value = '''0.07946'
$ws.Range('E34').Value = '  -3.02%  '
$ws.Range('D35').Value = '''0.9791'
$ws.Range('E35').Value = '  -6.73%  '
$ws.Range('D36').Value = '''0.02884'
$ws.Range('E36').Value = '  -4.92%  '
$ws.Range('D37').Value = '''6.572'
$ws.Range('E37').Value = '  -2.45%  '
$ws.Range('D38').Value = '''0.2648'
$ws.Range('E38').Value = '  -3.03%  '
$ws.Range('E39').Value = '  -7.93%  '
$ws.Range('D40').Value = '''0.09063'
$ws.Range('E40').Value = '  -1.34%  '
$ws.Range('D41').Value = '''0.7475'
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '''13.25'
$ws.Range('E42').Value = '  -3.81%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''1.416'
$ws.Range('E43').Value = '  -0.84%  '
$ws.Range('D44').Value = '''16.61'
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = '''0.6866'
$ws.Range('E45').Value = '  -2.75%  '
$ws.Range('D46').Value = '''2.402'
$ws.Range('E46').Value = '  -4.52%  '
$ws.Range('D47').Value = '''4.068'
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').Value = '''0.08202'
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').Value = '''133.87'
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('D51').Value = '''1.215'
$ws.Range('E51').Value = '  -2.45%  '
